# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.321.97"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "1.829.98"
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").Value = "'314.78"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").Value = "'0.4246"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("D8").Value = "'0.3697"
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("D9").Value = "'0.07274"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("D11").Value = "'21.13"
$ws.Range("E11").Value = "  -2.57%  "
$ws.Range("D12").Value = "1.827.22"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").Value = "'6.749"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.324"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.07093"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "'89.52"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "'1.006"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "'0.000008887"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").Value = "27.358.59"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").Value = "'5.145"
$ws.Range("E22").Value = "  -2.47%  "
$ws.Range("D23").Value = "'10.93"
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("D24").Value = "2.054.13"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").Value = "'1.997"
$ws.Range("D26").Value = "'153.04"
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").Value = "'2.183"
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("D28").Value = "'18.45"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").Value = "'5.255"
$ws.Range("E29").Value = "  -3.12%  "
$ws.Range("D30").Value = "'116.64"
$ws.Range("E30").Value = "  -3.63%  "
$ws.Range("D31").Value = "'0.08882"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("D32").Value = "'1.210"
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").Value = "'0.7613"
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("D34").Value = "'4.476"
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").Value = "'2.833"
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("D36").Value = "'1.005"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").Value = "'1.123"
$ws.Range("D38").Value = "'0.01984"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "'7.276"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D43").Value = "'0.5089"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").Value = "'8.710"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").Value = "'107.95"
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").Value = "'0.4780"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'1.005"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'0.06394"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("D50").Value = "'1.673"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("E51").Value = "  -3.20%  "
